$d = $word.ActiveDocument

# The preamble paragraph (2nd paragraph) holds both party-name/condition
# blocks that need to be updated:
#   {{ party1_name }}{% if party1_entity != "Individual" %}...
#   {{ party2_name }}{% if party2_entity != "Individual" %}...
# -> the bold/Helvetica-Neue "{{ partyN_name }}" run is folded into the
#    plain-styled run that follows it, and the Jinja condition gains a
#    .lower() call against a lower-cased literal.
$p = $d.Paragraphs.Item(2)
$pStart = $p.Range.Start

function Update-Party([string]$n) {
    $nameTok = "{{ party${n}_name }}"
    $oldCond = "{% if party${n}_entity != `"Individual`" %}, a {{ party${n}_state }} {{ party${n}_entity }}{% endif %}, whose address is {{ party${n}_address }}"
    $newCond = "{{ party${n}_name }}{% if party${n}_entity.lower() != `"individual`" %}, a {{ party${n}_state }} {{ party${n}_entity }}{% endif %}, whose address is {{ party${n}_address }}"

    # 1) Remove the bold "{{ partyN_name }}" run entirely.
    $text = $p.Range.Text
    $idx = $text.IndexOf($nameTok)
    $start = $pStart + $idx
    $end = $start + $nameTok.Length
    $d.Range($start, $end).Text = ""

    # 2) Rewrite the following (plain-styled) condition run so it starts
    #    with the party name and uses the new .lower() comparison. Doing
    #    this as a Range.Text assignment on the *existing* plain run keeps
    #    that run's own formatting (non-bold) instead of inheriting the
    #    bold formatting that used to precede it.
    $text2 = $p.Range.Text
    $idx2 = $text2.IndexOf($oldCond)
    $start2 = $pStart + $idx2
    $end2 = $start2 + $oldCond.Length
    $d.Range($start2, $end2).Text = $newCond
}

Update-Party "1"
Update-Party "2"
